$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 7-11 so only rows 1-6 remain
$ws.Range("A7:B11").EntireRow.Delete() | Out-Null

# Force text format on the data range so numeric-looking strings
# ("2", "1", "3", ...) are kept as text, not converted to numbers.
$dataRange = $ws.Range("A2:B6")
$dataRange.NumberFormat = "@"

# Update the data rows (A2:B6) with the new values
$ws.Range("A2").Value = "2"
$ws.Range("B2").Value = "1"

$ws.Range("A3").Value = "3"
$ws.Range("B3").Value = "4"

$ws.Range("A4").Value = "4"
$ws.Range("B4").Value = "3"

$ws.Range("A5").Value = "5"
$ws.Range("B5").Value = "2"

$ws.Range("A6").Value = "הדס"
$ws.Range("B6").Value = "6"
